$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161; this pushes the existing rows
# 161-180 down to 162-181 (matching the target dimension A1:R181).
$ws.Rows(161).Insert()

# Populate the newly inserted row 161 with the new weekly record.
$ws.Cells.Item(161, 1).Value2 = 1
$ws.Cells.Item(161, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(161, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(161, 4).Value2 = 44449
$ws.Cells.Item(161, 5).Value2 = 15
$ws.Cells.Item(161, 6).Value2 = 100114013
$ws.Cells.Item(161, 7).Value = "Zanahoria"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value2 = 96
$ws.Cells.Item(161, 11).Value2 = 8000
$ws.Cells.Item(161, 12).Value2 = 8500
$ws.Cells.Item(161, 13).Value2 = 8250
$ws.Cells.Item(161, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(161, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(161, 16).Value2 = 330
$ws.Cells.Item(161, 17).Value2 = 25
$ws.Cells.Item(161, 18).Value = "Hortaliza"
